$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "sudhir"
$ws.Range("A4").Value = "sachin"
$ws.Range("B3").Value = "127.0.0.2"
$ws.Range("B4").Value = "127.0.0.3"

$ws.Range("B3").Select()
